$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: Phone_Number -> text "792121074", entered with a leading apostrophe so
# it is kept as text (quote-prefixed) rather than being parsed as a number.
$ws.Range("B2").Value = "'792121074"

# A2: Name -> "MOHD1 Test23 Automation3 Fayoumi40"
$ws.Range("A2").Value = "MOHD1 Test23 Automation3 Fayoumi40"

# Move/save the selection at D9 (matches the author's final selection state)
$ws.Range("D9").Select()
